$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44508
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 10000
$ws.Range("P2").Value = 667

# Row 4
$ws.Range("D4").Value = 45225
$ws.Range("J4").Value = 80

# Row 5
$ws.Range("D5").Value = 44825
$ws.Range("J5").Value = 30

# Row 7
$ws.Range("D7").Value = 44757
$ws.Range("J7").Value = 30
$ws.Range("K7").Value = 20000
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = 20000
$ws.Range("P7").Value = 1333

# Row 8
$ws.Range("D8").Value = 44771
$ws.Range("J8").Value = 40
$ws.Range("K8").Value = 20000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 20000
$ws.Range("P8").Value = 1333

# Row 9
$ws.Range("D9").Value = 44749
$ws.Range("J9").Value = 50

# Row 10
$ws.Range("D10").Value = 44776
$ws.Range("J10").Value = 80

# Row 11
$ws.Range("D11").Value = 44819
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 20000
$ws.Range("L11").Value = 20000
$ws.Range("M11").Value = 20000
$ws.Range("P11").Value = 1333

# Row 12
$ws.Range("D12").Value = 44839
$ws.Range("J12").Value = 80
$ws.Range("K12").Value = 16000
$ws.Range("L12").Value = 16000
$ws.Range("M12").Value = 16000
$ws.Range("P12").Value = 1067

# Row 13
$ws.Range("D13").Value = 44767
$ws.Range("J13").Value = 50

# Row 14
$ws.Range("D14").Value = 44756
$ws.Range("J14").Value = 80

# Row 15
$ws.Range("D15").Value = 44812
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 20000
$ws.Range("L15").Value = 20000
$ws.Range("M15").Value = 20000
$ws.Range("P15").Value = 1333

# Row 16
$ws.Range("D16").Value = 44826
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = 20000
$ws.Range("L16").Value = 20000
$ws.Range("M16").Value = 20000
$ws.Range("P16").Value = 1333

# Row 17
$ws.Range("D17").Value = 44824
$ws.Range("J17").Value = 20

# Row 18
$ws.Range("D18").Value = 44827
$ws.Range("J18").Value = 20

# Row 19
$ws.Range("D19").Value = 44769
$ws.Range("J19").Value = 50

# Row 20
$ws.Range("D20").Value = 44518
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = 10000
$ws.Range("P20").Value = 667

# Row 21
$ws.Range("D21").Value = 44841
$ws.Range("J21").Value = 20
$ws.Range("K21").Value = 16000
$ws.Range("L21").Value = 16000
$ws.Range("M21").Value = 16000
$ws.Range("P21").Value = 1067

# Row 22
$ws.Range("D22").Value = 44811
$ws.Range("J22").Value = 30
$ws.Range("K22").Value = 20000
$ws.Range("L22").Value = 20000
$ws.Range("M22").Value = 20000
$ws.Range("P22").Value = 1333

# Row 23
$ws.Range("D23").Value = 44830
$ws.Range("J23").Value = 25
$ws.Range("K23").Value = 12000
$ws.Range("L23").Value = 12000
$ws.Range("M23").Value = 12000
$ws.Range("P23").Value = 800

# Row 24
$ws.Range("D24").Value = 45134
$ws.Range("J24").Value = 5

# Row 25
$ws.Range("D25").Value = 44755

# Row 26
$ws.Range("D26").Value = 44837
$ws.Range("J26").Value = 80
$ws.Range("K26").Value = 16000
$ws.Range("L26").Value = 16000
$ws.Range("M26").Value = 16000
$ws.Range("P26").Value = 1067

# Row 27
$ws.Range("D27").Value = 44525
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = 8000
$ws.Range("L27").Value = 8000
$ws.Range("M27").Value = 8000
$ws.Range("P27").Value = 533

# Row 28
$ws.Range("D28").Value = 44838
$ws.Range("J28").Value = 10
